$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 167, shifting existing rows 167-172 down to 168-173
$ws.Rows.Item(167).Insert()

# Populate new row 167 with the weekly data
$ws.Cells.Item(167, 1).Value2 = 11
$ws.Cells.Item(167, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(167, 3).Value2 = "Bíobío"
$ws.Cells.Item(167, 4).Value2 = 44706
$ws.Cells.Item(167, 5).Value2 = 8
$ws.Cells.Item(167, 6).Value2 = 100112003
$ws.Cells.Item(167, 7).Value2 = "Ajo"
$ws.Cells.Item(167, 8).Value2 = "Chino"
$ws.Cells.Item(167, 9).Value2 = "Primera"
$ws.Cells.Item(167, 10).Value2 = 400
$ws.Cells.Item(167, 11).Value2 = 17000
$ws.Cells.Item(167, 12).Value2 = 18000
$ws.Cells.Item(167, 13).Value2 = 17500
$ws.Cells.Item(167, 14).Value2 = "$/caja 10 kilos"
$ws.Cells.Item(167, 15).Value2 = "China"
$ws.Cells.Item(167, 16).Value2 = 1750
$ws.Cells.Item(167, 17).Value2 = 10
$ws.Cells.Item(167, 18).Value2 = "Hortaliza"

# Copy the date style (s="2") from the row below onto D167 to match formatting
$ws.Cells.Item(168, 4).Copy()
$ws.Cells.Item(167, 4).PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0
